$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.42
$ws.Range("G2").Value = 1.53
$ws.Range("H2").Value = 5.3
$ws.Range("J2").Value = 4.3
$ws.Range("K2").Value = 5.8
$ws.Range("P2").Value = 2.38
